# Applies bulk market-data refresh values across all profit sheets.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H19").Value = 1382.625
$ws.Range("I19").Value = 943
$ws.Range("J19").Value = 1646.4
$ws.Range("K19").Value = 943
$ws.Range("L19").Value = 1646.4
$ws.Range("M19").Value = -768
$ws.Range("N19").Value = -1996.4
$ws.Range("H28").Value = 1260.6
$ws.Range("I28").Value = 1293.9445
$ws.Range("K28").Value = 1293.9445
$ws.Range("M28").Value = -808.9445000000001
$ws.Range("H40").Value = 1295.46
$ws.Range("I40").Value = 1084.7273
$ws.Range("J40").Value = 1461.0358
$ws.Range("K40").Value = 1084.7273
$ws.Range("L40").Value = 1461.0358
$ws.Range("M40").Value = -909.7273
$ws.Range("N40").Value = -1811.0358
$ws.Range("H62").Value = 9614.541999999999
$ws.Range("I62").Value = 8989.444
$ws.Range("K62").Value = 8989.444
$ws.Range("M62").Value = -8365.444
$ws.Range("H65").Value = 9614.541999999999
$ws.Range("I65").Value = 8989.444
$ws.Range("K65").Value = 44947.22
$ws.Range("M65").Value = -41827.22
$ws.Range("H70").Value = 4087.25
$ws.Range("I70").Value = 2999.5
$ws.Range("J70").Value = 4449.8335
$ws.Range("K70").Value = 8998.5
$ws.Range("L70").Value = 13349.5005
$ws.Range("M70").Value = -8728.5
$ws.Range("N70").Value = -13889.5005
$ws.Range("H73").Value = 4087.25
$ws.Range("I73").Value = 2999.5
$ws.Range("J73").Value = 4449.8335
$ws.Range("K73").Value = 8998.5
$ws.Range("L73").Value = 13349.5005
$ws.Range("M73").Value = -8062.5
$ws.Range("N73").Value = -15221.5005
$ws.Range("H107").Value = 419.91177
$ws.Range("I107").Value = 419.91177
$ws.Range("K107").Value = 419.91177
$ws.Range("M107").Value = 1500.08823
$ws.Range("H127").Value = 500983.5
$ws.Range("I127").Value = 1750
$ws.Range("K127").Value = 5250
$ws.Range("M127").Value = -290

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 4221.3335
$ws.Range("I32").Value = 3382.2646
$ws.Range("K32").Value = 3382.2646
$ws.Range("M32").Value = -3095.2646
$ws.Range("H45").Value = 10529.833
$ws.Range("I45").Value = 17749.166
$ws.Range("K45").Value = 17749.166
$ws.Range("M45").Value = -17372.166
$ws.Range("H63").Value = 0
$ws.Range("I63").Value = 0
$ws.Range("K63").Value = 0
$ws.Range("M63").ClearContents()
$ws.Range("H66").Value = 0
$ws.Range("I66").Value = 0
$ws.Range("K66").Value = 0
$ws.Range("M66").ClearContents()
$ws.Range("H74").Value = 7917.3076
$ws.Range("I74").Value = 1441.5
$ws.Range("J74").Value = 18278.6
$ws.Range("K74").Value = 1441.5
$ws.Range("L74").Value = 18278.6
$ws.Range("M74").Value = -567.5
$ws.Range("N74").Value = -20026.6
$ws.Range("H76").Value = 0
$ws.Range("I76").Value = 0
$ws.Range("K76").Value = 0
$ws.Range("M76").ClearContents()
$ws.Range("H77").Value = 7917.3076
$ws.Range("I77").Value = 1441.5
$ws.Range("J77").Value = 18278.6
$ws.Range("K77").Value = 7207.5
$ws.Range("L77").Value = 91393
$ws.Range("M77").Value = -2839.5
$ws.Range("N77").Value = -100129
$ws.Range("H79").Value = 0
$ws.Range("I79").Value = 0
$ws.Range("K79").Value = 0
$ws.Range("M79").ClearContents()
$ws.Range("H88").Value = 3666.6667
$ws.Range("H91").Value = 3666.6667
$ws.Range("H122").Value = 64193.812
$ws.Range("I122").Value = 92053.91
$ws.Range("K122").Value = 276161.73
$ws.Range("M122").Value = -273711.73

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 3722.5483
$ws.Range("I20").Value = 3204.2917
$ws.Range("K20").Value = 3204.2917
$ws.Range("M20").Value = -2957.2917
$ws.Range("H134").Value = 2246.923
$ws.Range("I134").Value = 1536.4517
$ws.Range("J134").Value = 5000
$ws.Range("K134").Value = 4609.355100000001
$ws.Range("L134").Value = 15000
$ws.Range("M134").Value = -2074.355100000001
$ws.Range("N134").Value = -20070

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H58").Value = 3169.6365
$ws.Range("I58").Value = 3211.7144
$ws.Range("J58").Value = 3096
$ws.Range("K58").Value = 3211.7144
$ws.Range("L58").Value = 3096
$ws.Range("M58").Value = -3008.7144
$ws.Range("N58").Value = -3502
$ws.Range("H62").Value = 7667
$ws.Range("I62").Value = 8000.4
$ws.Range("J62").Value = 6000
$ws.Range("K62").Value = 8000.4
$ws.Range("L62").Value = 6000
$ws.Range("M62").Value = -7376.4
$ws.Range("N62").Value = -7248
$ws.Range("H65").Value = 7667
$ws.Range("I65").Value = 8000.4
$ws.Range("J65").Value = 6000
$ws.Range("K65").Value = 40002
$ws.Range("L65").Value = 30000
$ws.Range("M65").Value = -36882
$ws.Range("N65").Value = -36240
$ws.Range("H68").Value = 3068
$ws.Range("I68").Value = 3068
$ws.Range("K68").Value = 3068
$ws.Range("M68").Value = -2319
$ws.Range("H71").Value = 3068
$ws.Range("I71").Value = 3068
$ws.Range("K71").Value = 9204
$ws.Range("M71").Value = -5460
$ws.Range("H132").Value = 4093.6553
$ws.Range("I132").Value = 4106.143
$ws.Range("J132").Value = 4060.875
$ws.Range("K132").Value = 12318.429
$ws.Range("L132").Value = 12182.625
$ws.Range("M132").Value = -9788.429
$ws.Range("N132").Value = -17242.625
$ws.Range("H136").Value = 3169.6365
$ws.Range("I136").Value = 3211.7144
$ws.Range("J136").Value = 3096
$ws.Range("K136").Value = 9635.143199999999
$ws.Range("L136").Value = 9288
$ws.Range("M136").Value = -7085.143199999999
$ws.Range("N136").Value = -14388

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H97").Value = 565
$ws.Range("J97").Value = 525.1818
$ws.Range("L97").Value = 1575.5454
$ws.Range("N97").Value = -2567.5454
$ws.Range("H131").Value = 106800.81
$ws.Range("J131").Value = 17077
$ws.Range("L131").Value = 51231
$ws.Range("N131").Value = -61311

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H29").Value = 505000
$ws.Range("I29").Value = 505000
$ws.Range("K29").Value = 505000
$ws.Range("M29").Value = -504710
$ws.Range("H34").Value = 37499
$ws.Range("J34").Value = 37499
$ws.Range("L34").Value = 37499
$ws.Range("N34").Value = -38035
$ws.Range("H76").Value = 37499
$ws.Range("J76").Value = 37499
$ws.Range("L76").Value = 37499
$ws.Range("N76").Value = -38129
$ws.Range("H79").Value = 37499
$ws.Range("J79").Value = 37499
$ws.Range("L79").Value = 37499
$ws.Range("N79").Value = -39683
$ws.Range("H107").Value = 2967.318
$ws.Range("J107").Value = 3714.7144
$ws.Range("L107").Value = 3714.7144
$ws.Range("N107").Value = -7554.7144
$ws.Range("H113").Value = 5068.2856
$ws.Range("I113").Value = 3122
$ws.Range("K113").Value = 3122
$ws.Range("M113").Value = -952
$ws.Range("H122").Value = 3084
$ws.Range("I122").Value = 2523.3
$ws.Range("K122").Value = 7569.900000000001
$ws.Range("M122").Value = -5119.900000000001
$ws.Range("H126").Value = 21537.334
$ws.Range("I126").Value = 38431.668
$ws.Range("J126").Value = 4643
$ws.Range("K126").Value = 115295.004
$ws.Range("L126").Value = 13929
$ws.Range("M126").Value = -112825.004
$ws.Range("N126").Value = -18869

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 6465.9644
$ws.Range("I7").Value = 7264.143
$ws.Range("K7").Value = 7264.143
$ws.Range("M7").Value = -7152.143
$ws.Range("H16").Value = 16133639
$ws.Range("I16").Value = 41668040
$ws.Range("J16").Value = 6649.6313
$ws.Range("K16").Value = 41668040
$ws.Range("L16").Value = 6649.6313
$ws.Range("M16").Value = -41667870
$ws.Range("N16").Value = -6989.6313
$ws.Range("H19").Value = 4150
$ws.Range("I19").Value = 400
$ws.Range("J19").Value = 7900
$ws.Range("K19").Value = 400
$ws.Range("L19").Value = 7900
$ws.Range("M19").Value = -230
$ws.Range("N19").Value = -8240
$ws.Range("H40").Value = 5633.769
$ws.Range("I40").Value = 4822.9
$ws.Range("K40").Value = 4822.9
$ws.Range("M40").Value = -4686.9
$ws.Range("H48").Value = 0
$ws.Range("I48").Value = 0
$ws.Range("J48").Value = 0
$ws.Range("K48").Value = 0
$ws.Range("L48").Value = 0
$ws.Range("M48").ClearContents()
$ws.Range("N48").ClearContents()
$ws.Range("H68").Value = 17855
$ws.Range("I68").Value = 3972
$ws.Range("J68").Value = 27771.428
$ws.Range("K68").Value = 3972
$ws.Range("L68").Value = 27771.428
$ws.Range("M68").Value = -3223
$ws.Range("N68").Value = -29269.428
$ws.Range("H71").Value = 17855
$ws.Range("I71").Value = 3972
$ws.Range("J71").Value = 27771.428
$ws.Range("K71").Value = 19860
$ws.Range("L71").Value = 138857.14
$ws.Range("M71").Value = -16116
$ws.Range("N71").Value = -146345.14
$ws.Range("H126").Value = 6465.9644
$ws.Range("I126").Value = 7264.143
$ws.Range("K126").Value = 21792.429
$ws.Range("M126").Value = -19322.429

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H32").Value = 10000
$ws.Range("J32").Value = 0
$ws.Range("L32").Value = 0
$ws.Range("N32").ClearContents()
$ws.Range("H34").Value = 22200
$ws.Range("J34").Value = 22200
$ws.Range("L34").Value = 22200
$ws.Range("N34").Value = -22606
$ws.Range("H110").Value = 151000
$ws.Range("J110").Value = 151000
$ws.Range("L110").Value = 151000
$ws.Range("N110").Value = -159180
$ws.Range("H122").Value = 2371.516
$ws.Range("I122").Value = 2225.6072
$ws.Range("K122").Value = 6676.821599999999
$ws.Range("M122").Value = -4226.821599999999
$ws.Range("H132").Value = 2991.04
$ws.Range("I132").Value = 3120.6956
$ws.Range("K132").Value = 9362.086800000001
$ws.Range("M132").Value = -6832.086800000001
